# Insert a new data row at row 62, shifting the existing rows 62:162 down
# to 63:163 (matching columns/styles), then populate the new row with the
# new price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(62).Insert()

$ws.Range("A62").Value = 3
$ws.Range("B62").Value = "Femacal de La Calera"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 44540
$ws.Range("E62").Value = 5
$ws.Range("F62").Value = 100112010
$ws.Range("G62").Value = "Achicoria"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 80
$ws.Range("K62").Value = 5500
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = 5500
$ws.Range("N62").Value = "$/caja 16 unidades"
$ws.Range("O62").Value = "Provincia de Quillota"
$ws.Range("P62").Value = 344
$ws.Range("Q62").Value = 16
$ws.Range("R62").Value = "Hortaliza"
